# Re-update: strip the surrounding quotes from the DirectLayers / NodeLayers
# values (columns E and F) for each city row, and move the selection/scroll
# position back to a "fresh" view (A1 top-left, F8 selected) instead of the
# previously scrolled/selected AL4 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# City rows 3..7, columns E (DirectLayers) and F (NodeLayers):
# strip the double-quotes that used to wrap each comma-separated token.
# Write all of column E first, then all of column F, so newly-introduced
# shared-string entries land in the same append order the source file uses.
$ws.Cells.Item(3, 5).Value = "montreal_D_CD,montreal_D_CL"
$ws.Cells.Item(4, 5).Value = "vienna_D_CD,vienna_D_CL"
$ws.Cells.Item(5, 5).Value = "barcelona_D_CD,barcelona_D_CL,barcelona_D_EV"
$ws.Cells.Item(6, 5).Value = "budapest_D_CD,budapest_D_CL"
$ws.Cells.Item(7, 5).Value = "quebec_D_CD,quebec_D_CL"

$ws.Cells.Item(3, 6).Value = "montreal_N_CD,montreal_N_CL"
$ws.Cells.Item(4, 6).Value = "vienna_N_CD,vienna_N_CL"
$ws.Cells.Item(5, 6).Value = "barcelona_N_CD,barcelona_N_CL,barcelona_N_EV"
$ws.Cells.Item(6, 6).Value = "budapest_N_CD,budapest_N_CL"
$ws.Cells.Item(7, 6).Value = "quebec_N_CD,quebec_N_CL"

# Reset the view: scroll back so column A is visible again (was topLeftCell
# V1) and move the selection from AL4 to F8.
$ws.Range("A1").Select()
$ws.Range("F8").Select()
